# Updates the cryptos price/volume table with the latest scraped figures.
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Price cells that look like plain numbers are forced to Text format before
# the assignment (and the format/style reset afterwards) so values such as
# "244.07" or "1.00" are stored verbatim as strings instead of being
# re-interpreted as numeric doubles (which would lose trailing zeros /
# introduce float rounding noise). Two pairs of rows (7/8 and 46/47) have
# swapped rank order in the new data, so every field of those rows is
# rewritten explicitly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.448.45'
$ws.Range('E2').Value = '  -2.87%  '
$ws.Range('D3').Value = '1.972.90'
$ws.Range('E3').Value = '  -4.07%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.07'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.633'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.46%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.48'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.05'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.356'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0729'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.91%  '
$ws.Range('E12').Value = '  -2.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.941'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.20'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.62%  '
$ws.Range('D15').Value = '2.263.17'
$ws.Range('E15').Value = '  -3.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.24'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.62%  '
$ws.Range('D17').Value = '1.958.12'
$ws.Range('E17').Value = '  -4.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.43'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.29%  '
$ws.Range('D19').Value = '35.382.50'
$ws.Range('E19').Value = '  -2.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.33'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('D21').Value = '0.0₃0838'
$ws.Range('E21').Value = '  -2.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '231.87'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.10'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.56%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.53'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +19.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.28'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.49'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.02'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.02'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.52%  '
$ws.Range('E30').Value = '  -3.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.82'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.10'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0587'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('E34').Value = '  +11.39%  '
$ws.Range('E35').Value = '  -3.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.35'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.24%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  -3.71%  '
$ws.Range('E39').Value = '  +5.27%  '
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('E41').Value = '  +1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0209'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.52%  '
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '90.61'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.34%  '
$ws.Range('D45').Value = '1.370.70'
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0879'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.71%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.73'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.45'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.81'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.31%  '
$ws.Range('E51').Value = '  +11.10%  '
